# Rename the "Contrasts" sheet to "peripostinterval" and repoint the
# workbook's defined names at the new sheet name (mirrors opening the
# workbook in a newer Excel, renaming the sheet, and re-saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contrasts")
$ws.Name = "peripostinterval"

# Renaming the sheet automatically repoints formula-style defined names
# (e.g. "peripostinterval" -> peripostinterval!$A$1:$I$7), but the ones
# that refer to #REF! lose their sheet-qualifier on rename, so restore it
# explicitly to match the target workbook.
$wb.Names.Item("aucContr").RefersTo = "=peripostinterval!#REF!"
$wb.Names.Item("aucDiffContr").RefersTo = "=peripostinterval!#REF!"
$wb.Names.Item("periGroup").RefersTo = "=peripostinterval!#REF!"
$wb.Names.Item("perigroupdiff").RefersTo = "=peripostinterval!#REF!"
